# Auto-generated Excel COM-interop script
# Applies the scheduled-runner data refresh described in the commit diff:
# updates cached currentAveragePrice / LevePrice* / LeveProfit* columns (H,I,J,K,L,M,N)
# for the affected Leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(94, 8).Value = 5533.5  # H94: 4577.8 -> 5533.5
$ws.Cells.Item(94, 9).Value = 5533.5  # I94: 4577.8 -> 5533.5
$ws.Cells.Item(94, 11).Value = 5533.5  # K94: 4577.8 -> 5533.5
$ws.Cells.Item(94, 13).Value = -5082.5  # M94: -4126.8 -> -5082.5

$ws.Cells.Item(116, 8).Value = 8156.2  # H116: 8008.4116 -> 8156.2
$ws.Cells.Item(116, 9).Value = 7081.8887  # I116: 7048.8184 -> 7081.8887
$ws.Cells.Item(116, 11).Value = 7081.8887  # K116: 7048.8184 -> 7081.8887
$ws.Cells.Item(116, 13).Value = -3639.8887  # M116: -3606.8184 -> -3639.8887

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5042.4614  # H2: 5413.0835 -> 5042.4614
$ws.Cells.Item(2, 9).Value = 1493.8462  # I2: 1567.5 -> 1493.8462
$ws.Cells.Item(2, 10).Value = 8591.076999999999  # J2: 9258.666999999999 -> 8591.076999999999
$ws.Cells.Item(2, 11).Value = 1493.8462  # K2: 1567.5 -> 1493.8462
$ws.Cells.Item(2, 12).Value = 8591.076999999999  # L2: 9258.666999999999 -> 8591.076999999999
$ws.Cells.Item(2, 13).Value = -1380.8462  # M2: -1454.5 -> -1380.8462
$ws.Cells.Item(2, 14).Value = -8817.076999999999  # N2: -9484.666999999999 -> -8817.076999999999

$ws.Cells.Item(32, 8).Value = 4198.052  # H32: 4225.2026 -> 4198.052
$ws.Cells.Item(32, 9).Value = 3880.137  # I32: 3981.1973 -> 3880.137
$ws.Cells.Item(32, 11).Value = 3880.137  # K32: 3981.1973 -> 3880.137
$ws.Cells.Item(32, 13).Value = -3593.137  # M32: -3694.1973 -> -3593.137

$ws.Cells.Item(45, 8).Value = 3620.5  # H45: 3854.875 -> 3620.5
$ws.Cells.Item(45, 9).Value = 3288.25  # I45: 3451.5454 -> 3288.25
$ws.Cells.Item(45, 10).Value = 4285  # J45: 4742.2 -> 4285
$ws.Cells.Item(45, 11).Value = 3288.25  # K45: 3451.5454 -> 3288.25
$ws.Cells.Item(45, 12).Value = 4285  # L45: 4742.2 -> 4285
$ws.Cells.Item(45, 13).Value = -2911.25  # M45: -3074.5454 -> -2911.25
$ws.Cells.Item(45, 14).Value = -5039  # N45: -5496.2 -> -5039

$ws.Cells.Item(97, 8).Value = 864.7931  # H97: 1053.0322 -> 864.7931
$ws.Cells.Item(97, 9).Value = 824.5833  # I97: 894.2 -> 824.5833
$ws.Cells.Item(97, 10).Value = 1057.8  # J97: 1714.8334 -> 1057.8
$ws.Cells.Item(97, 11).Value = 824.5833  # K97: 894.2 -> 824.5833
$ws.Cells.Item(97, 12).Value = 1057.8  # L97: 1714.8334 -> 1057.8
$ws.Cells.Item(97, 13).Value = -328.5833  # M97: -398.2 -> -328.5833
$ws.Cells.Item(97, 14).Value = -2049.8  # N97: -2706.8334 -> -2049.8

$ws.Cells.Item(116, 8).Value = 5042.4614  # H116: 5413.0835 -> 5042.4614
$ws.Cells.Item(116, 9).Value = 1493.8462  # I116: 1567.5 -> 1493.8462
$ws.Cells.Item(116, 10).Value = 8591.076999999999  # J116: 9258.666999999999 -> 8591.076999999999
$ws.Cells.Item(116, 11).Value = 1493.8462  # K116: 1567.5 -> 1493.8462
$ws.Cells.Item(116, 12).Value = 8591.076999999999  # L116: 9258.666999999999 -> 8591.076999999999
$ws.Cells.Item(116, 13).Value = 800.1538  # M116: 726.5 -> 800.1538
$ws.Cells.Item(116, 14).Value = -13179.077  # N116: -13846.667 -> -13179.077

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5042.4614  # H3: 5413.0835 -> 5042.4614
$ws.Cells.Item(3, 9).Value = 1493.8462  # I3: 1567.5 -> 1493.8462
$ws.Cells.Item(3, 10).Value = 8591.076999999999  # J3: 9258.666999999999 -> 8591.076999999999
$ws.Cells.Item(3, 11).Value = 1493.8462  # K3: 1567.5 -> 1493.8462
$ws.Cells.Item(3, 12).Value = 8591.076999999999  # L3: 9258.666999999999 -> 8591.076999999999
$ws.Cells.Item(3, 13).Value = -1379.8462  # M3: -1453.5 -> -1379.8462
$ws.Cells.Item(3, 14).Value = -8819.076999999999  # N3: -9486.666999999999 -> -8819.076999999999

$ws.Cells.Item(9, 8).Value = 0  # H9: 29998 -> 0
$ws.Cells.Item(9, 10).Value = 0  # J9: 29998 -> 0
$ws.Cells.Item(9, 12).Value = 0  # L9: 29998 -> 0
$ws.Cells.Item(9, 14).ClearContents()  # N9: -30334 -> (removed)

$ws.Cells.Item(134, 8).Value = 3869.25  # H134: 3972.8572 -> 3869.25
$ws.Cells.Item(134, 9).Value = 3639.1538  # I134: 3729.182 -> 3639.1538
$ws.Cells.Item(134, 11).Value = 10917.4614  # K134: 11187.546 -> 10917.4614
$ws.Cells.Item(134, 13).Value = -8382.4614  # M134: -8652.545999999998 -> -8382.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 30000  # H3: 24999.75 -> 30000
$ws.Cells.Item(3, 9).Value = 0  # I3: 20000 -> 0
$ws.Cells.Item(3, 10).Value = 30000  # J3: 26666.334 -> 30000
$ws.Cells.Item(3, 11).Value = 0  # K3: 20000 -> 0
$ws.Cells.Item(3, 12).Value = 30000  # L3: 26666.334 -> 30000
$ws.Cells.Item(3, 13).ClearContents()  # M3: -19887 -> (removed)
$ws.Cells.Item(3, 14).Value = -30226  # N3: -26892.334 -> -30226

$ws.Cells.Item(31, 8).Value = 37995.965  # H31: 37989.133 -> 37995.965
$ws.Cells.Item(31, 9).Value = 2355.0557  # I31: 2422.7058 -> 2355.0557
$ws.Cells.Item(31, 10).Value = 91457.336  # J31: 84499.08 -> 91457.336
$ws.Cells.Item(31, 11).Value = 2355.0557  # K31: 2422.7058 -> 2355.0557
$ws.Cells.Item(31, 12).Value = 91457.336  # L31: 84499.08 -> 91457.336
$ws.Cells.Item(31, 13).Value = -2060.0557  # M31: -2127.7058 -> -2060.0557
$ws.Cells.Item(31, 14).Value = -92047.336  # N31: -85089.08 -> -92047.336

$ws.Cells.Item(34, 8).Value = 37995.965  # H34: 37989.133 -> 37995.965
$ws.Cells.Item(34, 9).Value = 2355.0557  # I34: 2422.7058 -> 2355.0557
$ws.Cells.Item(34, 10).Value = 91457.336  # J34: 84499.08 -> 91457.336
$ws.Cells.Item(34, 11).Value = 2355.0557  # K34: 2422.7058 -> 2355.0557
$ws.Cells.Item(34, 12).Value = 91457.336  # L34: 84499.08 -> 91457.336
$ws.Cells.Item(34, 13).Value = -2153.0557  # M34: -2220.7058 -> -2153.0557
$ws.Cells.Item(34, 14).Value = -91861.336  # N34: -84903.08 -> -91861.336

$ws.Cells.Item(58, 8).Value = 4101.154  # H58: 4110.1284 -> 4101.154
$ws.Cells.Item(58, 9).Value = 2517.7407  # I58: 2530.7036 -> 2517.7407
$ws.Cells.Item(58, 11).Value = 2517.7407  # K58: 2530.7036 -> 2517.7407
$ws.Cells.Item(58, 13).Value = -2314.7407  # M58: -2327.7036 -> -2314.7407

$ws.Cells.Item(74, 8).Value = 40000  # H74: 0 -> 40000
$ws.Cells.Item(74, 10).Value = 40000  # J74: 0 -> 40000
$ws.Cells.Item(74, 12).Value = 40000  # L74: 0 -> 40000
$ws.Cells.Item(74, 14).Value = -41748  # N74: (new) -> -41748

$ws.Cells.Item(77, 8).Value = 40000  # H77: 0 -> 40000
$ws.Cells.Item(77, 10).Value = 40000  # J77: 0 -> 40000
$ws.Cells.Item(77, 12).Value = 120000  # L77: 0 -> 120000
$ws.Cells.Item(77, 14).Value = -128736  # N77: (new) -> -128736

$ws.Cells.Item(88, 8).Value = 0  # H88: 17499.5 -> 0
$ws.Cells.Item(88, 10).Value = 0  # J88: 17499.5 -> 0
$ws.Cells.Item(88, 12).Value = 0  # L88: 17499.5 -> 0
$ws.Cells.Item(88, 14).ClearContents()  # N88: -18311.5 -> (removed)

$ws.Cells.Item(91, 8).Value = 0  # H91: 17499.5 -> 0
$ws.Cells.Item(91, 10).Value = 0  # J91: 17499.5 -> 0
$ws.Cells.Item(91, 12).Value = 0  # L91: 17499.5 -> 0
$ws.Cells.Item(91, 14).ClearContents()  # N91: -20307.5 -> (removed)

$ws.Cells.Item(132, 8).Value = 4124.346  # H132: 4172.4 -> 4124.346
$ws.Cells.Item(132, 9).Value = 3425.7917  # I132: 3601.7896 -> 3425.7917
$ws.Cells.Item(132, 10).Value = 12507  # J132: 15014 -> 12507
$ws.Cells.Item(132, 11).Value = 10277.3751  # K132: 10805.3688 -> 10277.3751
$ws.Cells.Item(132, 12).Value = 37521  # L132: 45042 -> 37521
$ws.Cells.Item(132, 13).Value = -7747.375100000001  # M132: -8275.3688 -> -7747.375100000001
$ws.Cells.Item(132, 14).Value = -42581  # N132: -50102 -> -42581

$ws.Cells.Item(136, 8).Value = 4101.154  # H136: 4110.1284 -> 4101.154
$ws.Cells.Item(136, 9).Value = 2517.7407  # I136: 2530.7036 -> 2517.7407
$ws.Cells.Item(136, 11).Value = 7553.222099999999  # K136: 7592.110799999999 -> 7553.222099999999
$ws.Cells.Item(136, 13).Value = -5003.222099999999  # M136: -5042.110799999999 -> -5003.222099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 68  # H33: 77.40000000000001 -> 68
$ws.Cells.Item(33, 9).Value = 49.75  # I33: 59.333332 -> 49.75
$ws.Cells.Item(33, 11).Value = 298.5  # K33: 355.999992 -> 298.5
$ws.Cells.Item(33, 13).Value = -15.5  # M33: -72.99999200000002 -> -15.5

$ws.Cells.Item(51, 8).Value = 5  # H51: 2000 -> 5
$ws.Cells.Item(51, 9).Value = 0  # I51: 2000 -> 0
$ws.Cells.Item(51, 10).Value = 5  # J51: 0 -> 5
$ws.Cells.Item(51, 11).Value = 0  # K51: 6000 -> 0
$ws.Cells.Item(51, 12).Value = 15  # L51: 0 -> 15
$ws.Cells.Item(51, 13).ClearContents()  # M51: -5540 -> (removed)
$ws.Cells.Item(51, 14).Value = -935  # N51: (new) -> -935

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 295228.6  # H3: 295232.47 -> 295228.6
$ws.Cells.Item(3, 10).Value = 1470.6364  # J3: 1476.6364 -> 1470.6364
$ws.Cells.Item(3, 12).Value = 1470.6364  # L3: 1476.6364 -> 1470.6364
$ws.Cells.Item(3, 14).Value = -1702.6364  # N3: -1708.6364 -> -1702.6364

$ws.Cells.Item(20, 8).Value = 46152.11  # H20: 46365.75 -> 46152.11
$ws.Cells.Item(20, 10).Value = 51110.4  # J20: 52777.25 -> 51110.4
$ws.Cells.Item(20, 12).Value = 51110.4  # L20: 52777.25 -> 51110.4
$ws.Cells.Item(20, 14).Value = -51600.4  # N20: -53267.25 -> -51600.4

$ws.Cells.Item(80, 8).Value = 1256748  # H80: 559934.9 -> 1256748
$ws.Cells.Item(80, 9).Value = 5000000  # I80: 2500900 -> 5000000
$ws.Cells.Item(80, 10).Value = 8997.333000000001  # J80: 5373.4287 -> 8997.333000000001
$ws.Cells.Item(80, 11).Value = 5000000  # K80: 2500900 -> 5000000
$ws.Cells.Item(80, 12).Value = 8997.333000000001  # L80: 5373.4287 -> 8997.333000000001
$ws.Cells.Item(80, 13).Value = -4999002  # M80: -2499902 -> -4999002
$ws.Cells.Item(80, 14).Value = -10993.333  # N80: -7369.4287 -> -10993.333

$ws.Cells.Item(83, 8).Value = 1256748  # H83: 559934.9 -> 1256748
$ws.Cells.Item(83, 9).Value = 5000000  # I83: 2500900 -> 5000000
$ws.Cells.Item(83, 10).Value = 8997.333000000001  # J83: 5373.4287 -> 8997.333000000001
$ws.Cells.Item(83, 11).Value = 25000000  # K83: 12504500 -> 25000000
$ws.Cells.Item(83, 12).Value = 44986.665  # L83: 26867.1435 -> 44986.665
$ws.Cells.Item(83, 13).Value = -24995008  # M83: -12499508 -> -24995008
$ws.Cells.Item(83, 14).Value = -54970.665  # N83: -36851.14350000001 -> -54970.665

$ws.Cells.Item(113, 8).Value = 3724.35  # H113: 4084.8823 -> 3724.35
$ws.Cells.Item(113, 9).Value = 2823.625  # I113: 3465.6667 -> 2823.625
$ws.Cells.Item(113, 10).Value = 4324.8335  # J113: 4422.636 -> 4324.8335
$ws.Cells.Item(113, 11).Value = 2823.625  # K113: 3465.6667 -> 2823.625
$ws.Cells.Item(113, 12).Value = 4324.8335  # L113: 4422.636 -> 4324.8335
$ws.Cells.Item(113, 13).Value = -653.625  # M113: -1295.6667 -> -653.625
$ws.Cells.Item(113, 14).Value = -8664.833500000001  # N113: -8762.636 -> -8664.833500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 4263.222  # H46: 3564.25 -> 4263.222
$ws.Cells.Item(46, 9).Value = 2092  # I46: 1933.6 -> 2092
$ws.Cells.Item(46, 10).Value = 6000.2  # J46: 4729 -> 6000.2
$ws.Cells.Item(46, 11).Value = 2092  # K46: 1933.6 -> 2092
$ws.Cells.Item(46, 12).Value = 6000.2  # L46: 4729 -> 6000.2
$ws.Cells.Item(46, 13).Value = -1904  # M46: -1745.6 -> -1904
$ws.Cells.Item(46, 14).Value = -6376.2  # N46: -5105 -> -6376.2

$ws.Cells.Item(132, 8).Value = 6078.4546  # H132: 5961.4414 -> 6078.4546
$ws.Cells.Item(132, 9).Value = 4959.4443  # I132: 4808.9473 -> 4959.4443
$ws.Cells.Item(132, 11).Value = 14878.3329  # K132: 14426.8419 -> 14878.3329
$ws.Cells.Item(132, 13).Value = -12348.3329  # M132: -11896.8419 -> -12348.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(26, 8).Value = 8500  # H26: 0 -> 8500
$ws.Cells.Item(26, 10).Value = 8500  # J26: 0 -> 8500
$ws.Cells.Item(26, 12).Value = 8500  # L26: 0 -> 8500
$ws.Cells.Item(26, 14).Value = -9086  # N26: (new) -> -9086

$ws.Cells.Item(126, 8).Value = 2544.258  # H126: 2450.6177 -> 2544.258
$ws.Cells.Item(126, 9).Value = 1601.96  # I126: 1613.4231 -> 1601.96
$ws.Cells.Item(126, 10).Value = 6470.5  # J126: 5171.5 -> 6470.5
$ws.Cells.Item(126, 11).Value = 4805.88  # K126: 4840.2693 -> 4805.88
$ws.Cells.Item(126, 12).Value = 19411.5  # L126: 15514.5 -> 19411.5
$ws.Cells.Item(126, 13).Value = -2335.88  # M126: -2370.2693 -> -2335.88
$ws.Cells.Item(126, 14).Value = -24351.5  # N126: -20454.5 -> -24351.5

$ws.Cells.Item(136, 8).Value = 2529.257  # H136: 2596.4412 -> 2529.257
$ws.Cells.Item(136, 9).Value = 1103.8334  # I136: 1133.4482 -> 1103.8334
$ws.Cells.Item(136, 11).Value = 3311.5002  # K136: 3400.3446 -> 3311.5002
$ws.Cells.Item(136, 13).Value = -761.5001999999999  # M136: -850.3446000000004 -> -761.5001999999999
